$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.183.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.76%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.849.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.67%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'232.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.90%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = '  +0.12%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4674"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.85%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = '  -4.79%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.06369"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'1.880.22"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.96%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07417"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.79%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'16.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.85%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'4.937"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.53%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'84.90"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.74%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.6269"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.25%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'30.128.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.85%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = '  +0.11%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'228.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.46%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'12.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.42%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.000007308"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.10%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'2.099.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.41%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.15%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'4.936"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.00%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'5.916"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.14%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'9.220"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.68%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'165.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.40%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = '  -4.99%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'1.864"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.73%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = '  +2.93%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.387"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.89%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'4.106"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.82%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = '  -4.19%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.04882"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.74%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.156"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.53%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.7122"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = "'1.001"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.691"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.77%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.01851"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.83%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'2.633"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.03%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.9071"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.28%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'1.942"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.94%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'105.05"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.77%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.9985"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.60%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'5.536"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.05%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.4066"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.47%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'7.013"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'60.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.80%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.1185"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -7.04%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'8.575"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.62%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'32.93"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.11%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.05571"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.65%  '
$ws.Range("E51").Style = "Normal"
